$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C2:C7) from 2023-10-22 (45221) to 2023-10-25 (45224)
$ws.Range("C2:C7").Value = 45224
